$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting rows 8-13 down to 9-14
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 content
$ws.Range("A8").Value = "被"
$ws.Range("B8").Value = "passive"
$ws.Range("C8").Value = "虛詞"

# A8 should carry the same formatting as B9 (old B8 "left+top+bottom" red box style)
$ws.Range("B9").Copy()
$ws.Range("A8").PasteSpecial(-4122)

# C8 should carry the same formatting as C9 (old C8 "right+top+bottom" red box style)
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# B8 needs a new combined style: red font (matching A8/C8 font) + top/bottom-only
# medium red border (no left/right), completing the boxed row.
$rng = $ws.Range("B8")
$rng.Font.Color = 255
$rng.Font.Charset = 136
$rng.Font.Family = 1
$rng.Borders.Item(8).Weight = -4138
$rng.Borders.Item(8).Color = 255
$rng.Borders.Item(9).Weight = -4138
$rng.Borders.Item(9).Color = 255

# Match the new selection state: A8:C8, active cell A8
[void]$ws.Range("A8:C8").Select()

Write-Host "done"
